$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$edits = @(
    @("H17", 621),
    @("I17", 0),
    @("J17", 621),
    @("K17", 0),
    @("L17", 1863),
    @("M17", "__DELETE__"),
    @("N17", -2199),
    @("H21", 30759.777),
    @("J21", 21133.334),
    @("L21", 21133.334),
    @("N21", -22069.334),
    @("H23", 30759.777),
    @("J23", 21133.334),
    @("L23", 21133.334),
    @("N23", -21601.334),
    @("H43", 857.4286),
    @("I43", 797.7143),
    @("J43", 887.2857),
    @("K43", 797.7143),
    @("L43", 887.2857),
    @("M43", -728.7143),
    @("N43", -1025.2857),
    @("H64", 3703.3),
    @("I64", 3564.7058),
    @("J64", 3884.5386),
    @("K64", 3564.7058),
    @("L64", 3884.5386),
    @("M64", -3316.7058),
    @("N64", -4380.5386),
    @("H67", 3703.3),
    @("I67", 3564.7058),
    @("J67", 3884.5386),
    @("K67", 3564.7058),
    @("L67", 3884.5386),
    @("M67", -2706.7058),
    @("N67", -5600.5386),
    @("H76", 3666),
    @("I76", 3000),
    @("J76", 3999),
    @("K76", 3000),
    @("L76", 3999),
    @("M76", -2685),
    @("N76", -4629),
    @("H79", 3666),
    @("I79", 3000),
    @("J79", 3999),
    @("K79", 3000),
    @("L79", 3999),
    @("M79", -1908),
    @("N79", -6183),
    @("H98", 1051),
    @("I98", 984.3684),
    @("J98", 1367.5),
    @("K98", 984.3684),
    @("L98", 1367.5),
    @("M98", 513.6316),
    @("N98", -4363.5),
    @("H112", 1344),
    @("J112", 1403.129),
    @("L112", 4209.387),
    @("N112", -6425.387),
    @("H115", 10000668),
    @("I115", 10000668),
    @("K115", 30002004),
    @("M115", -30000437),
    @("H116", 3080595.2),
    @("I116", 19233020),
    @("J116", 3942.8572),
    @("K116", 19233020),
    @("L116", 3942.8572),
    @("M116", -19229578),
    @("N116", -10826.8572),
    @("H122", 1051),
    @("I122", 984.3684),
    @("J122", 1367.5),
    @("K122", 2953.1052),
    @("L122", 4102.5),
    @("M122", -503.1052),
    @("N122", -9002.5),
    @("H125", 3241.6667),
    @("I125", 2254.3333),
    @("J125", 3570.7778),
    @("K125", 20288.9997),
    @("L125", 32137.0002),
    @("M125", -17828.9997),
    @("N125", -37057.00019999999),
    @("H132", 2286.4333),
    @("I132", 1840.9412),
    @("J132", 4810.8887),
    @("K132", 5522.8236),
    @("L132", 14432.6661),
    @("M132", -2992.8236),
    @("N132", -19492.6661),
    @("H135", 880.05554),
    @("I135", 614.0625),
    @("J135", 3008),
    @("K135", 5526.5625),
    @("L135", 27072),
    @("M135", -2991.5625),
    @("N135", -32142),
    @("H140", 36320),
    @("J140", 36320),
    @("L140", 36320),
    @("N140", -46680),
    @("H141", 2293),
    @("I141", 1610.2),
    @("J141", 4000),
    @("K141", 4830.6),
    @("L141", 12000),
    @("M141", 349.3999999999996),
    @("N141", -22360)
)
foreach ($pair in $edits) {
    $ref = $pair[0]
    $val = $pair[1]
    if ("$val" -eq "__DELETE__") {
        $ws.Range($ref).ClearContents()
    } else {
        $ws.Range($ref).Value = $val
    }
}

$ws = $wb.Worksheets.Item("ARM")
$edits = @(
    @("H62", 5226),
    @("I62", 5226),
    @("J62", 0),
    @("K62", 5226),
    @("L62", 0),
    @("N62", "__DELETE__"),
    @("M62", -4602),
    @("H65", 5226),
    @("I65", 5226),
    @("J65", 0),
    @("K65", 15678),
    @("L65", 0),
    @("N65", "__DELETE__"),
    @("M65", -12558)
)
foreach ($pair in $edits) {
    $ref = $pair[0]
    $val = $pair[1]
    if ("$val" -eq "__DELETE__") {
        $ws.Range($ref).ClearContents()
    } else {
        $ws.Range($ref).Value = $val
    }
}

$ws = $wb.Worksheets.Item("BSM")
$edits = @(
    @("H99", 2128.6),
    @("I99", 1484.2858),
    @("J99", 2692.375),
    @("K99", 1484.2858),
    @("L99", 2692.375),
    @("M99", 13.71419999999989),
    @("N99", -5688.375)
)
foreach ($pair in $edits) {
    $ref = $pair[0]
    $val = $pair[1]
    if ("$val" -eq "__DELETE__") {
        $ws.Range($ref).ClearContents()
    } else {
        $ws.Range($ref).Value = $val
    }
}

$ws = $wb.Worksheets.Item("CRP")
$edits = @(
    @("H58", 935.12),
    @("I58", 760.93335),
    @("K58", 760.93335),
    @("M58", -557.93335),
    @("H136", 935.12),
    @("I136", 760.93335),
    @("K136", 2282.80005),
    @("M136", 267.1999500000002)
)
foreach ($pair in $edits) {
    $ref = $pair[0]
    $val = $pair[1]
    if ("$val" -eq "__DELETE__") {
        $ws.Range($ref).ClearContents()
    } else {
        $ws.Range($ref).Value = $val
    }
}

$ws = $wb.Worksheets.Item("CUL")
$edits = @(
    @("H3", 4127.5),
    @("I3", 2947.1428),
    @("K3", 8841.4284),
    @("M3", -8729.4284),
    @("H103", 525),
    @("I103", 525),
    @("J103", 0),
    @("K103", 1575),
    @("L103", 0),
    @("N103", "__DELETE__"),
    @("M103", -696),
    @("H131", 644.8788),
    @("I131", 423.72223),
    @("J131", 910.26666),
    @("K131", 1271.16669),
    @("L131", 2730.79998),
    @("M131", 3768.83331),
    @("N131", -12810.79998),
    @("H137", 3167.0667),
    @("I137", 2983.4285),
    @("J137", 3327.75),
    @("K137", 8950.2855),
    @("L137", 9983.25),
    @("M137", -3850.2855),
    @("N137", -20183.25)
)
foreach ($pair in $edits) {
    $ref = $pair[0]
    $val = $pair[1]
    if ("$val" -eq "__DELETE__") {
        $ws.Range($ref).ClearContents()
    } else {
        $ws.Range($ref).Value = $val
    }
}

$ws = $wb.Worksheets.Item("GSM")
$edits = @(
    @("H102", 1143.1333),
    @("I102", 1134.3846),
    @("J102", 1200),
    @("K102", 1134.3846),
    @("L102", 1200),
    @("M102", 487.6153999999999),
    @("N102", -4444),
    @("H122", 2500),
    @("I122", 2000),
    @("J122", 3000),
    @("K122", 6000),
    @("L122", 9000),
    @("M122", -3550),
    @("N122", -13900)
)
foreach ($pair in $edits) {
    $ref = $pair[0]
    $val = $pair[1]
    if ("$val" -eq "__DELETE__") {
        $ws.Range($ref).ClearContents()
    } else {
        $ws.Range($ref).Value = $val
    }
}

$ws = $wb.Worksheets.Item("LTW")
$edits = @(
    @("H22", 974),
    @("I22", 974),
    @("J22", 0),
    @("K22", 974),
    @("L22", 0),
    @("M22", -679),
    @("N22", "__DELETE__"),
    @("H26", 2009),
    @("I26", 2009),
    @("K26", 2009),
    @("M26", -1714),
    @("H27", 974),
    @("I27", 974),
    @("J27", 0),
    @("K27", 974),
    @("L27", 0),
    @("M27", -867),
    @("N27", "__DELETE__"),
    @("H32", 0),
    @("I32", 0),
    @("J32", 0),
    @("K32", 0),
    @("L32", 0),
    @("M32", "__DELETE__"),
    @("N32", "__DELETE__"),
    @("H33", 1500),
    @("J33", 0),
    @("L33", 0),
    @("N33", "__DELETE__"),
    @("H34", 12000),
    @("J34", 12000),
    @("L34", 12000),
    @("N34", -12344),
    @("H40", 8030.3),
    @("I40", 10351),
    @("J40", 6483.1665),
    @("K40", 10351),
    @("L40", 6483.1665),
    @("M40", -10215),
    @("N40", -6755.1665),
    @("H122", 3569.3845),
    @("I122", 3460.3333),
    @("J122", 3662.8572),
    @("K122", 10380.9999),
    @("L122", 10988.5716),
    @("M122", -7930.999899999999),
    @("N122", -15888.5716)
)
foreach ($pair in $edits) {
    $ref = $pair[0]
    $val = $pair[1]
    if ("$val" -eq "__DELETE__") {
        $ws.Range($ref).ClearContents()
    } else {
        $ws.Range($ref).Value = $val
    }
}

$ws = $wb.Worksheets.Item("WVR")
$edits = @(
    @("H46", 0),
    @("J46", 0),
    @("L46", 0),
    @("N46", "__DELETE__"),
    @("H134", 0),
    @("J134", 0),
    @("L134", 0),
    @("N134", "__DELETE__")
)
foreach ($pair in $edits) {
    $ref = $pair[0]
    $val = $pair[1]
    if ("$val" -eq "__DELETE__") {
        $ws.Range($ref).ClearContents()
    } else {
        $ws.Range($ref).Value = $val
    }
}
